# Update countries & provincias Spain
# - Re-order "Armenia" to sit right after "Moldavia" (row 64), shifting
#   Ghana -> row 65 and Finlandia -> row 66, and refresh Armenia's stats
#   with up-to-date numbers. Ghana/Finlandia keep the numbers that used to
#   belong to the row below them (the feed just re-sorted + refreshed).
# - Bump the "Datos actualizados" timestamp in A1 from 08:35 to 09:05.
# - Refresh Letonia's (row 108) stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 09:05"

# Row 64: Armenia (new/updated figures)
$ws.Range("A64").Value = "Armenia"
$ws.Range("B64").Value = 6661
$ws.Range("C64").Value = 359
$ws.Range("D64").Value = 3064
$ws.Range("E64").Value = 3516
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 4
$ws.Range("H64").Value = 81

# Row 65: Ghana (shifted down, keeps the numbers formerly on row 64)
$ws.Range("A65").Value = "Ghana"
$ws.Range("B65").Value = 6617
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 1978
$ws.Range("E65").Value = 4608
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 31

# Row 66: Finlandia (shifted down, keeps the numbers formerly on row 65)
$ws.Range("A66").Value = "Finlandia"
$ws.Range("B66").Value = 6568
$ws.Range("C66").Value = 0
$ws.Range("D66").Value = 4800
$ws.Range("E66").Value = 1462
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 306

# Row 108: Letonia (updated figures)
$ws.Range("B108").Value = 1047
$ws.Range("C108").Value = 1
$ws.Range("E108").Value = 313
